$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 6.368
$ws.Range("A12").Value = -21.435
$ws.Range("D14").Value = -8.158000000000001
$ws.Range("D19").Value = -7.994
$ws.Range("B23").Value = 7.628
$ws.Range("D24").Value = -7.379
$ws.Range("A27").Value = -21.637
$ws.Range("B28").Value = 5.36
$ws.Range("A32").Value = -21.137
$ws.Range("B32").Value = 6.955
$ws.Range("B34").Value = 6.425
$ws.Range("A36").Value = -20.764
$ws.Range("A38").Value = -20.722
$ws.Range("D38").Value = -8.392999999999999
$ws.Range("D41").Value = -8.209
$ws.Range("B42").Value = 8.046000000000001
$ws.Range("A46").Value = -21.578
$ws.Range("B49").Value = 6.255999999999999
$ws.Range("D52").Value = -7.941000000000001
$ws.Range("A54").Value = -20.886
$ws.Range("B54").Value = 5.731999999999999
$ws.Range("A55").Value = -22.016
$ws.Range("A56").Value = -21.644
$ws.Range("A67").Value = -21.422
$ws.Range("A69").Value = -21.387
$ws.Range("A72").Value = -21.57
$ws.Range("D72").Value = -7.734
$ws.Range("B78").Value = 6.921000000000001
$ws.Range("D78").Value = -8.116
$ws.Range("B80").Value = 7.519
$ws.Range("A83").Value = -21.216
$ws.Range("D83").Value = -7.978999999999999
$ws.Range("D85").Value = -8.643000000000001
$ws.Range("A86").Value = -21.733
$ws.Range("D86").Value = -8.497
$ws.Range("D90").Value = -7.045
$ws.Range("A91").Value = -20.94
$ws.Range("A93").Value = -21.545
$ws.Range("D96").Value = -7.279999999999999
$ws.Range("B97").Value = 5.491000000000001
$ws.Range("A99").Value = -20.858
$ws.Range("B99").Value = 5.944
$ws.Range("B101").Value = 5.252000000000001
$ws.Range("D103").Value = -8.342000000000002
$ws.Range("A104").Value = -21.437